$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates ---------------------------------------------
# The "Special Effects" labels are being renamed from the mEFCT_SPCL|...
# naming scheme to mEFCT_SPCLFUNC|...
$ws.Range("D4").Value = "mEFCT_SPCLFUNC|mEFCT_SHOOT"
$ws.Range("E4").Value = "mEFCT_SPCLFUNC|mEFCT_SHOOT"

$ws.Range("D5").Value = "mEFCT_SPCLFUNC|mEFCT_OPEN_BARREL"
$ws.Range("E5").Value = "mEFCT_SPCLFUNC|mEFCT_OPEN_BARREL"

$ws.Range("D6").Value = "mEFCT_SPCLFUNC|mEFCT_LOCK_LOAD"
$ws.Range("E6").Value = "mEFCT_SPCLFUNC|mEFCT_LOCK_LOAD"

# --- Column widths -------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 44.592447916666664
$ws.Columns("E:E").ColumnWidth = 37.451822916666664

# --- Row 24 height reverts to default (was manually set to 30) -----------
$ws.Rows("24:24").AutoFit() | Out-Null

# --- Selection / active cell ---------------------------------------------
$ws.Range("C6").Select()
